$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a 8-row x 20-column array holding the new data for rows 2..9 (columns A..T)
$data = New-Object 'object[,]' 8,20

$data[0,0] = "ECs"
$data[0,1] = "Glg1"
$data[0,2] = "Sele"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 21.62921866666666
$data[0,7] = 64.88765599999999
$data[0,8] = 0.1903297306418182
$data[0,9] = 0.1903297306418182
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 4.699506666666667
$data[0,13] = 14.09852
$data[0,14] = 0.9660495246229048
$data[0,15] = 0.9660495246229047
$data[0,16] = 101.6466573187911
$data[0,17] = 914.8199158691199
$data[0,18] = 0.183867945808134
$data[0,19] = 0.183867945808134

$data[1,0] = "ECs"
$data[1,1] = "Glg1"
$data[1,2] = "Sele"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 21.62921866666666
$data[1,7] = 64.88765599999999
$data[1,8] = 0.1903297306418182
$data[1,9] = 0.1903297306418182
$data[1,10] = 1
$data[1,11] = 0.3333333333333333
$data[1,12] = 0.1651576666666667
$data[1,13] = 0.495473
$data[1,14] = 0.03395047537709522
$data[1,15] = 0.03395047537709522
$data[1,16] = 3.572231286809777
$data[1,17] = 32.150081581288
$data[1,18] = 0.006461784833684215
$data[1,19] = 0.006461784833684214

$data[2,0] = "FAPs"
$data[2,1] = "Glg1"
$data[2,2] = "Sele"
$data[2,3] = "ECs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 53.56207000000001
$data[2,7] = 160.68621
$data[2,8] = 0.4713279066076088
$data[2,9] = 0.4713279066076088
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 4.699506666666667
$data[2,13] = 14.09852
$data[2,14] = 0.9660495246229048
$data[2,15] = 0.9660495246229047
$data[2,16] = 251.7153050454667
$data[2,17] = 2265.4377454092
$data[2,18] = 0.4553261001197894
$data[2,19] = 0.4553261001197893

$data[3,0] = "FAPs"
$data[3,1] = "Glg1"
$data[3,2] = "Sele"
$data[3,3] = "FAPs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 53.56207000000001
$data[3,7] = 160.68621
$data[3,8] = 0.4713279066076088
$data[3,9] = 0.4713279066076088
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.1651576666666667
$data[3,13] = 0.495473
$data[3,14] = 0.03395047537709522
$data[3,15] = 0.03395047537709522
$data[3,16] = 8.846186503036668
$data[3,17] = 79.61567852733
$data[3,18] = 0.01600180648781946
$data[3,19] = 0.01600180648781946

$data[4,0] = "M2"
$data[4,1] = "Glg1"
$data[4,2] = "Sele"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 21.64753366666666
$data[4,7] = 64.942601
$data[4,8] = 0.1904908963811095
$data[4,9] = 0.1904908963811095
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 4.699506666666667
$data[4,13] = 14.09852
$data[4,14] = 0.9660495246229048
$data[4,15] = 0.9660495246229047
$data[4,16] = 101.7327287833911
$data[4,17] = 915.5945590505199
$data[4,18] = 0.1840236398939618
$data[4,19] = 0.1840236398939618

$data[5,0] = "M2"
$data[5,1] = "Glg1"
$data[5,2] = "Sele"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 21.64753366666666
$data[5,7] = 64.942601
$data[5,8] = 0.1904908963811095
$data[5,9] = 0.1904908963811095
$data[5,10] = 1
$data[5,11] = 0.3333333333333333
$data[5,12] = 0.1651576666666667
$data[5,13] = 0.495473
$data[5,14] = 0.03395047537709522
$data[5,15] = 0.03395047537709522
$data[5,16] = 3.575256149474777
$data[5,17] = 32.177305345273
$data[5,18] = 0.006467256487147653
$data[5,19] = 0.006467256487147653

$data[6,0] = "sCs"
$data[6,1] = "Glg1"
$data[6,2] = "Sele"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 16.80195566666667
$data[6,7] = 50.405867
$data[6,8] = 0.1478514663694635
$data[6,9] = 0.1478514663694635
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 4.699506666666667
$data[6,13] = 14.09852
$data[6,14] = 0.9660495246229048
$data[6,15] = 0.9660495246229047
$data[6,16] = 78.9609026685378
$data[6,17] = 710.64812401684
$data[6,18] = 0.1428318388010196
$data[6,19] = 0.1428318388010196

$data[7,0] = "sCs"
$data[7,1] = "Glg1"
$data[7,2] = "Sele"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 16.80195566666667
$data[7,7] = 50.405867
$data[7,8] = 0.1478514663694635
$data[7,9] = 0.1478514663694635
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.1651576666666667
$data[7,13] = 0.495473
$data[7,14] = 0.03395047537709522
$data[7,15] = 0.03395047537709522
$data[7,16] = 2.774971793343445
$data[7,17] = 24.974746140091
$data[7,18] = 0.005019627568443892
$data[7,19] = 0.005019627568443891

$ws.Range("A2:T9").Value = $data
